$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update K column (최종점수) values
$ws.Range("K2").Value = 54.9
$ws.Range("K3").Value = 50.9
$ws.Range("K4").Value = 46.7
$ws.Range("K5").Value = 45.5
$ws.Range("K6").Value = 36.9

# Update N column (MACRO_SCORE) values
$ws.Range("N2").Value = 51.53902399942638
$ws.Range("N3").Value = 51.53902399942638
$ws.Range("N4").Value = 51.53902399942638
$ws.Range("N5").Value = 51.53902399942638
$ws.Range("N6").Value = 51.53902399942638
